$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns containing numeric-looking text must be forced to Text format
# so Excel does not auto-convert them to actual numbers.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "42.666.56"
$ws.Range("E2").Value = "  -0.14%  "
$ws.Range("D3").Value = "2.540.66"
$ws.Range("E3").Value = "  +0.68%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "309.07"
$ws.Range("E5").Value = "  -2.26%  "
$ws.Range("D6").Value = "97.38"
$ws.Range("E6").Value = "  +0.35%  "
$ws.Range("D7").Value = "0.571"
$ws.Range("E7").Value = "  -0.74%  "
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("E9").Value = "  -0.61%  "
$ws.Range("D10").Value = "35.55"
$ws.Range("E10").Value = "  -0.61%  "
$ws.Range("E11").Value = "  -0.65%  "
$ws.Range("D12").Value = "7.39"
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").Value = "0.108"
$ws.Range("E13").Value = "  -2.14%  "
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "2.932.12"
$ws.Range("E14").Value = "  +0.61%  "
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "2.611.24"
$ws.Range("E15").Value = "  +3.34%  "
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").Value = "15.65"
$ws.Range("E16").Value = "  +3.84%  "
$ws.Range("D17").Value = "0.833"
$ws.Range("E17").Value = "  -1.48%  "
$ws.Range("D18").Value = "42.641.69"
$ws.Range("E18").Value = "  -0.39%  "
$ws.Range("D19").Value = "6.72"
$ws.Range("E19").Value = "  -1.30%  "
$ws.Range("B20").Value = "InternetComputer(DFINITY)"
$ws.Range("C20").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D20").Value = "12.38"
$ws.Range("E20").Value = "  -2.82%  "
$ws.Range("B21").Value = "ShibaInu"
$ws.Range("C21").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D21").Value = "0.0₃0954"
$ws.Range("E21").Value = "  -0.59%  "
$ws.Range("D22").Value = "69.22"
$ws.Range("E22").Value = "  -0.57%  "
$ws.Range("D23").Value = "246.76"
$ws.Range("E23").Value = "  -1.85%  "
$ws.Range("E24").Value = "  -1.31%  "
$ws.Range("E25").Value = "  -0.26%  "
$ws.Range("D26").Value = "26.49"
$ws.Range("E26").Value = "  +0.57%  "
$ws.Range("E27").Value = "  +0.02%  "
$ws.Range("E28").Value = "  -1.50%  "
$ws.Range("D29").Value = "39.91"
$ws.Range("E29").Value = "  -2.44%  "
$ws.Range("D30").Value = "10.14"
$ws.Range("E30").Value = "  -2.00%  "
$ws.Range("D31").Value = "157.63"
$ws.Range("E31").Value = "  -1.31%  "
$ws.Range("D32").Value = "5.71"
$ws.Range("E32").Value = "  -3.35%  "
$ws.Range("D33").Value = "0.0792"
$ws.Range("E33").Value = "  +0.49%  "
$ws.Range("D34").Value = "3.29"
$ws.Range("E34").Value = "  -0.70%  "
$ws.Range("E35").Value = "  -3.71%  "
$ws.Range("D36").Value = "2.61"
$ws.Range("E36").Value = "  -4.05%  "
$ws.Range("D37").Value = "18.42"
$ws.Range("E37").Value = "  -2.12%  "
$ws.Range("E38").Value = "  +7.43%  "
$ws.Range("E39").Value = "  -1.50%  "
$ws.Range("E40").Value = "  -0.50%  "
$ws.Range("D41").Value = "22.48"
$ws.Range("E41").Value = "  +3.21%  "
$ws.Range("D42").Value = "4.06"
$ws.Range("E42").Value = "  +6.47%  "
$ws.Range("E43").Value = "  -0.18%  "
$ws.Range("E44").Value = "  -1.86%  "
$ws.Range("D45").Value = "1.985.72"
$ws.Range("E45").Value = "  -1.91%  "
$ws.Range("D46").Value = "3.19"
$ws.Range("E46").Value = "  -1.71%  "
$ws.Range("D47").Value = "8.94"
$ws.Range("E47").Value = "  -1.12%  "
$ws.Range("D48").Value = "2.797.04"
$ws.Range("E48").Value = "  +1.00%  "
$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D49").Value = "0.193"
$ws.Range("E49").Value = "  +1.63%  "
$ws.Range("B50").Value = "BitcoinSV"
$ws.Range("C50").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D50").Value = "80.61"
$ws.Range("E50").Value = "  -4.06%  "
$ws.Range("D51").Value = "73.37"
$ws.Range("E51").Value = "  -2.18%  "

# Remove the temporary text-number-format so no stray cell styles remain
$ws.Range("D4").ClearFormats()
$ws.Range("D5").ClearFormats()
$ws.Range("D6").ClearFormats()
$ws.Range("D7").ClearFormats()
$ws.Range("D10").ClearFormats()
$ws.Range("D12").ClearFormats()
$ws.Range("D13").ClearFormats()
$ws.Range("D16").ClearFormats()
$ws.Range("D17").ClearFormats()
$ws.Range("D19").ClearFormats()
$ws.Range("D20").ClearFormats()
$ws.Range("D22").ClearFormats()
$ws.Range("D23").ClearFormats()
$ws.Range("D26").ClearFormats()
$ws.Range("D29").ClearFormats()
$ws.Range("D30").ClearFormats()
$ws.Range("D31").ClearFormats()
$ws.Range("D32").ClearFormats()
$ws.Range("D33").ClearFormats()
$ws.Range("D34").ClearFormats()
$ws.Range("D36").ClearFormats()
$ws.Range("D37").ClearFormats()
$ws.Range("D41").ClearFormats()
$ws.Range("D42").ClearFormats()
$ws.Range("D46").ClearFormats()
$ws.Range("D47").ClearFormats()
$ws.Range("D49").ClearFormats()
$ws.Range("D50").ClearFormats()
$ws.Range("D51").ClearFormats()
